$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "30÷3=10, 0"
$t.Cell(1,2).Range.Text = "44÷4=11, 0"
$t.Cell(1,3).Range.Text = "60÷7=8, 4"
$t.Cell(1,4).Range.Text = "79÷5=15, 4"
$t.Cell(1,5).Range.Text = "56÷5=11, 1"

$t.Cell(5,1).Range.Text = "98÷5=19, 3"
$t.Cell(5,2).Range.Text = "84÷3=28, 0"
$t.Cell(5,3).Range.Text = "97÷2=48, 1"
$t.Cell(5,4).Range.Text = "70÷2=35, 0"
$t.Cell(5,5).Range.Text = "80÷8=10, 0"

$t.Cell(9,1).Range.Text = "81÷2=40, 1"
$t.Cell(9,2).Range.Text = "48÷9=5, 3"
$t.Cell(9,3).Range.Text = "33÷8=4, 1"
$t.Cell(9,4).Range.Text = "66÷2=33, 0"
$t.Cell(9,5).Range.Text = "35÷3=11, 2"

$t.Cell(13,1).Range.Text = "57÷2=28, 1"
$t.Cell(13,2).Range.Text = "20÷7=2, 6"
$t.Cell(13,3).Range.Text = "51÷6=8, 3"
$t.Cell(13,4).Range.Text = "40÷9=4, 4"
$t.Cell(13,5).Range.Text = "93÷5=18, 3"

$t.Cell(17,1).Range.Text = "80÷2=40, 0"
$t.Cell(17,2).Range.Text = "34÷5=6, 4"
$t.Cell(17,3).Range.Text = "24÷9=2, 6"
$t.Cell(17,4).Range.Text = "36÷5=7, 1"
$t.Cell(17,5).Range.Text = "96÷6=16, 0"

